# Update cryptocurrency price (D) and volume change (E) columns
# to reflect refreshed data, matching the upstream commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.897.66'
$ws.Range('E2').Value = '  -0.32%  '
$ws.Range('D3').Value = '2.118.13'
$ws.Range('E3').Value = '  +0.81%  '
$ws.Range('D4').Value = "'" + '1.006'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').Value = "'" + '348.29'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.00%  '
$ws.Range('D6').Value = "'" + '1.006'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.40%  '
$ws.Range('D7').Value = "'" + '0.5191'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.04%  '
$ws.Range('D8').Value = "'" + '0.4464'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.90%  '
$ws.Range('D9').Value = "'" + '54.13'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.75%  '
$ws.Range('D10').Value = "'" + '0.09379'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.51%  '
$ws.Range('D11').Value = "'" + '1.181'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.29%  '
$ws.Range('D12').Value = "'" + '25.21'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.31%  '
$ws.Range('D13').Value = '2.142.68'
$ws.Range('E13').Value = '  +2.29%  '
$ws.Range('D14').Value = "'" + '8.420'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.91%  '
$ws.Range('D15').Value = "'" + '6.858'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.09%  '
$ws.Range('D16').Value = "'" + '102.50'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.83%  '
$ws.Range('D17').Value = "'" + '0.00001166'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.89%  '
$ws.Range('D18').Value = "'" + '1.007'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.29%  '
$ws.Range('D19').Value = "'" + '21.58'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.20%  '
$ws.Range('D20').Value = "'" + '0.06676'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.13%  '
$ws.Range('D21').Value = "'" + '6.308'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.53%  '
$ws.Range('E22').Value = '  +0.45%  '
$ws.Range('D23').Value = '29.925.34'
$ws.Range('E23').Value = '  -0.56%  '
$ws.Range('D24').Value = "'" + '12.72'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.22%  '
$ws.Range('D25').Value = "'" + '2.328'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.17%  '
$ws.Range('D26').Value = '2.352.25'
$ws.Range('E26').Value = '  +0.34%  '
$ws.Range('D27').Value = "'" + '22.13'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.01%  '
$ws.Range('D28').Value = "'" + '2.570'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.35%  '
$ws.Range('D29').Value = "'" + '162.63'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.30%  '
$ws.Range('D30').Value = "'" + '134.10'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.81%  '
$ws.Range('D31').Value = "'" + '1.157'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.12%  '
$ws.Range('D32').Value = "'" + '1.804'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +9.70%  '
$ws.Range('D33').Value = "'" + '0.1056'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.12%  '
$ws.Range('D34').Value = "'" + '6.261'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.72%  '
$ws.Range('D35').Value = "'" + '3.976'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('D36').Value = "'" + '6.517'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.37%  '
$ws.Range('D37').Value = "'" + '10.83'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +7.30%  '
$ws.Range('D38').Value = "'" + '0.02604'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.47%  '
$ws.Range('D39').Value = "'" + '0.06824'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.74%  '
$ws.Range('D40').Value = "'" + '12.70'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.89%  '
$ws.Range('D41').Value = "'" + '0.7029'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.33%  '
$ws.Range('D42').Value = "'" + '1.346'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.99%  '
$ws.Range('D43').Value = "'" + '0.2248'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.78%  '
$ws.Range('D44').Value = "'" + '0.6850'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.77%  '
$ws.Range('D45').Value = "'" + '14.49'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.61%  '
$ws.Range('D46').Value = "'" + '2.359'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.59%  '
$ws.Range('E47').Value = '  +0.49%  '
$ws.Range('E48').Value = '  +0.07%  '
$ws.Range('D49').Value = "'" + '3.635'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.19%  '
$ws.Range('D50').Value = "'" + '1.253'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.29%  '
$ws.Range('E51').Value = '  +0.68%  '
